# Gemeentes.xlsx - add an "Inhabitants" column (E) to Blad1.
#
# 1. Add the new shared string "Inhabitants " and use it as the header
#    for column E (E1).
# 2. Fill in the inhabitant counts for every data row (E2:E71). The value
#    only depends on the (region, year) pair - it repeats for the two
#    vehicle-type rows ("Snorfietsen"/"Bromfietsen") that share a year.
# 3. A few cells (E13:E15) carried a stray "Helvetica 12pt" style left
#    over from earlier edits; rows 65:71 never had an E cell at all.
#    Both are fixed by copying D's (already-correct) format into E
#    before writing the value, so every E cell ends up sharing the same
#    style as the rest of the column.
# 4. Leave the selection where the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. header -------------------------------------------------------
$ws.Range("E1").Value = "Inhabitants "

# --- 3. normalise formatting for the rows that need a format fix-up --
# E13:E15 had an orphaned style; E65:E71 had no cell at all yet.
$fixupRows = @(13, 14, 15, 65, 66, 67, 68, 69, 70, 71)
foreach ($r in $fixupRows) {
    $ws.Range("D" + $r).Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)
}

# --- 2. inhabitant counts, keyed by row -------------------------------
$inhabitants = @{
    2=57715;    3=58055;    4=58524;    5=58846;    6=60341;    7=60370;    8=61003;
    9=57715;   10=58055;   11=58524;   12=58846;   13=60341;   14=60370;   15=61003;
   16=537833;  17=545838;  18=548320;  19=553417;  20=562839;  21=566221;  22=568945;
   23=537833;  24=545838;  25=548320;  26=553417;  27=562839;  28=566221;  29=568945;
   30=123107;  31=124084;  32=124481;  33=125504;  34=127073;  35=128810;  36=129973;
   37=123107;  38=124084;  39=124481;  40=125504;  41=127073;  42=128810;  43=129973;
   44=176731;  45=177659;  46=177359;  47=179100;  48=182480;  49=187049;  50=189007;
   51=176731;  52=177659;  53=177359;  54=179100;  55=182480;  56=187049;  57=189007;
   58=644618;  59=651157;  60=651631;  61=655468;  62=663900;  63=670610;  64=672960;
   65=644618;  66=651157;  67=651631;  68=655468;  69=663900;  70=670610;  71=672960
}

foreach ($r in $inhabitants.Keys) {
    $ws.Cells.Item($r, 5).Value = $inhabitants[$r]
}

# --- 4. selection left by the author on save --------------------------
$ws.Range("G68").Select()
